# Update "想去人数" (interest count) figures in the "展览" and "全部类型"
# sheets to match the newly generated site output.
#
# 展览 (sheet 1):
#   F2  266  -> 267
#   F15 1305 -> 1306
#   F24 5603 -> 5607
#   F29 14130 -> 14135
#   F35 573  -> 576
#   F37 105  -> 106
#   F39 110  -> 111
#
# 全部类型 (sheet 4) has the same events (duplicated rows / offset by the
# extra rows from other categories) so it needs the identical value bumps:
#   F2  266  -> 267
#   F15 1305 -> 1306
#   F27 5603 -> 5607
#   F32 14130 -> 14135
#   F38 573  -> 576
#   F40 105  -> 106
#   F42 110  -> 111

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 267
$ws1.Range("F15").Value = 1306
$ws1.Range("F24").Value = 5607
$ws1.Range("F29").Value = 14135
$ws1.Range("F35").Value = 576
$ws1.Range("F37").Value = 106
$ws1.Range("F39").Value = 111

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 267
$ws4.Range("F15").Value = 1306
$ws4.Range("F27").Value = 5607
$ws4.Range("F32").Value = 14135
$ws4.Range("F38").Value = 576
$ws4.Range("F40").Value = 106
$ws4.Range("F42").Value = 111
